$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old "Comment" column (F), shifting it to H.
$ws.Columns("F:G").Insert()

# New data row (14) - comment text is added first so the new shared string
# table matches the author's original insertion order.
$ws.Range("H14").Value = "frequency"

# New header cells for the inserted columns.
$ws.Range("F1").Value = "clipp"
$ws.Range("G1").Value = "dec"

# New "clipp" numeric column values for the existing data rows (2-13).
$ws.Range("F2:F9").Value = 5
$ws.Range("F10:F13").Value = 50
$ws.Range("F2:F13").HorizontalAlignment = -4108

# New "dec" text column values for the existing data rows (2-13).
$ws.Range("G2:G13").Value = "e-5"
$ws.Range("G2:G13").ClearFormats()

# Remaining values for the new data row (14).
$ws.Range("B14").Value = 4
$ws.Range("C14").Value = 6
$ws.Range("D14").Value = 1000
$ws.Range("E14").Value = 0.29
$ws.Range("F14").Value = 15
$ws.Range("G14").Value = "e-5"
$ws.Range("B14:F14").HorizontalAlignment = -4108

# Match the author's final selection.
$ws.Range("F15").Select() | Out-Null
